$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# Row 5
$ws.Range("M5").Value = 1.07
$ws.Range("O5").Value = 1.47

# Row 6
$ws.Range("M6").Value = 1.11
$ws.Range("O6").Value = 1.58

# Row 7
$ws.Range("M7").Value = 1.08
$ws.Range("O7").Value = 1.47

# Row 8
$ws.Range("M8").Value = 1.13
$ws.Range("O8").Value = 1.58

# Row 9
$ws.Range("BD9").Value = 126
$ws.Range("M9").Value = 1.05
$ws.Range("O9").Value = 1.41
$ws.Range("P9").Value = 2.62

# Row 10
$ws.Range("M10").Value = 1.04
$ws.Range("O10").Value = 1.27

# Row 18
$ws.Range("U18").Value = 1.72

# Row 29
$ws.Range("AB29").Value = 29
$ws.Range("AH29").Value = 26
$ws.Range("AM29").Value = 51
$ws.Range("G29").Value = 1.22
$ws.Range("I29").Value = 13

# Row 30
$ws.Range("AH30").Value = 15
$ws.Range("AI30").Value = 17
$ws.Range("AL30").Value = 17
$ws.Range("G30").Value = 2.5
$ws.Range("I30").Value = 2.55
$ws.Range("Y30").Value = 11

# Row 32
$ws.Range("AC32").Value = 7
$ws.Range("AH32").Value = 8
$ws.Range("G32").Value = 1.95
$ws.Range("H32").Value = 3.5
$ws.Range("I32").Value = 3.8
$ws.Range("J32").Value = 2.75
$ws.Range("M32").Value = 1.1
$ws.Range("N32").Value = 7
$ws.Range("Q32").Value = 2.5
$ws.Range("R32").Value = 1.5
$ws.Range("X32").Value = 8
$ws.Range("Z32").Value = 17

# Row 34
$ws.Range("G34").Value = 1.17
$ws.Range("U34").Value = 1.92
$ws.Range("V34").Value = 1.77

# Row 35
$ws.Range("I35").Value = 1.33
$ws.Range("U35").Value = 1.69

# Row 36
$ws.Range("I36").Value = 1.62
$ws.Range("U36").Value = 1.47

# Row 37
$ws.Range("G37").Value = 2.1
$ws.Range("Q37").Value = 1.94
$ws.Range("R37").Value = 1.79
$ws.Range("U37").Value = 1.77
$ws.Range("V37").Value = 1.87

# Row 38
$ws.Range("G38").Value = 1.71
$ws.Range("U38").Value = 1.5
$ws.Range("V38").Value = 2.37

# Row 39
$ws.Range("AC39").Value = 11
$ws.Range("AE39").Value = 21
$ws.Range("AF39").Value = 67
$ws.Range("AH39").Value = 6.5
$ws.Range("AI39").Value = 6.5
$ws.Range("AJ39").Value = 9
$ws.Range("AK39").Value = 9
$ws.Range("AM39").Value = 29
$ws.Range("AS39").Value = 351
$ws.Range("AT39").Value = 3
$ws.Range("AU39").Value = 9.5
$ws.Range("AV39").Value = 67
$ws.Range("AW39").Value = 3.25
$ws.Range("AX39").Value = 6.5
$ws.Range("BB39").Value = 151
$ws.Range("G39").Value = 8
$ws.Range("I39").Value = 1.4
$ws.Range("J39").Value = 7.5
$ws.Range("K39").Value = 2.38
$ws.Range("M39").Value = 1.05
$ws.Range("N39").Value = 11
$ws.Range("O39").Value = 1.25
$ws.Range("P39").Value = 3.75
$ws.Range("Q39").Value = 1.8
$ws.Range("R39").Value = 2
$ws.Range("S39").Value = 1.36
$ws.Range("T39").Value = 3
$ws.Range("U39").Value = 2.1
$ws.Range("V39").Value = 1.67
$ws.Range("W39").Value = 19
$ws.Range("Z39").Value = 101

# Row 42
$ws.Range("M42").Value = 1.07
$ws.Range("N42").Value = 9
$ws.Range("Q42").Value = 2.25
$ws.Range("R42").Value = 1.62

# Row 44
$ws.Range("U44").Value = 1.77
$ws.Range("V44").Value = 1.87

# Row 45
$ws.Range("U45").Value = 1.87
$ws.Range("V45").Value = 1.77

# Row 46
$ws.Range("AT46").Value = 2.75
$ws.Range("O46").Value = 1.3
$ws.Range("P46").Value = 3.4
$ws.Range("R46").Value = 1.8
$ws.Range("S46").Value = 1.4
$ws.Range("T46").Value = 2.75
$ws.Range("U46").Value = 1.8
$ws.Range("V46").Value = 1.8

# Row 47
$ws.Range("U47").Value = 1.69

# Row 55
$ws.Range("M55").Value = 1.08
$ws.Range("N55").Value = 8

# Row 56
$ws.Range("M56").Value = 1.04
$ws.Range("O56").Value = 1.2

# Row 57
$ws.Range("M57").Value = 1.04
$ws.Range("O57").Value = 1.22

# Row 58
$ws.Range("J58").Value = 2.62
$ws.Range("M58").Value = 1.02
$ws.Range("O58").Value = 1.13

# Row 59
$ws.Range("M59").Value = 1.07
$ws.Range("O59").Value = 1.33

# Row 60
$ws.Range("M60").Value = 1.03
$ws.Range("O60").Value = 1.18

# Row 64
$ws.Range("U64").Value = 1.54

# Row 65
$ws.Range("U65").Value = 1.54

# Row 66
$ws.Range("M66").Value = 1.03
$ws.Range("O66").Value = 1.19

# Row 67
$ws.Range("M67").Value = 1.01
$ws.Range("O67").Value = 1.08

# Row 68
$ws.Range("M68").Value = 1.02
$ws.Range("O68").Value = 1.13

# Row 69
$ws.Range("AW69").Value = 4.75
$ws.Range("M69").Value = 1.02
$ws.Range("O69").Value = 1.13

# Row 70
$ws.Range("M70").Value = 1.03
$ws.Range("O70").Value = 1.17

# Row 71
$ws.Range("M71").Value = 1.02
$ws.Range("O71").Value = 1.13

# Row 72
$ws.Range("M72").Value = 1.05
$ws.Range("O72").Value = 1.33

# Row 73
$ws.Range("M73").Value = 1.05
$ws.Range("O73").Value = 1.33

# Row 74
$ws.Range("V74").Value = 1.69

# Row 75
$ws.Range("U75").Value = 1.8
$ws.Range("V75").Value = 1.8

# Row 76
$ws.Range("V76").Value = 1.63

# Row 77
$ws.Range("M77").Value = 1.04
$ws.Range("O77").Value = 1.22
$ws.Range("U77").Value = 1.72

# Row 81
$ws.Range("M81").Value = 1.08
$ws.Range("O81").Value = 1.4
$ws.Range("R81").Value = 1.58

# Row 82
$ws.Range("M82").Value = 1.1
$ws.Range("O82").Value = 1.5
$ws.Range("R82").Value = 1.44

# Row 83
$ws.Range("M83").Value = 1.04
$ws.Range("O83").Value = 1.22
$ws.Range("U83").Value = 1.8
$ws.Range("V83").Value = 1.8

# Row 84
$ws.Range("V84").Value = 1.63

# Row 85
$ws.Range("U85").Value = 1.63

# Row 88
$ws.Range("Q88").Value = 1.17

# Row 90
$ws.Range("Q90").Value = 1.77

# Row 91
$ws.Range("R91").Value = 1.63

# Row 92
$ws.Range("Q92").Value = 1.5

# Row 94
$ws.Range("AA94").Value = 15.5
$ws.Range("AB94").Value = 20
$ws.Range("AC94").Value = 8.75
$ws.Range("AD94").Value = 7.2
$ws.Range("AE94").Value = 11.25
$ws.Range("AH94").Value = 13.5
$ws.Range("AI94").Value = 19
$ws.Range("AJ94").Value = 10.5
$ws.Range("AK94").Value = 40
$ws.Range("AL94").Value = 22
$ws.Range("AM94").Value = 23
$ws.Range("AN94").Value = 4.4
$ws.Range("AO94").Value = 11
$ws.Range("AP94").Value = 16.5
$ws.Range("AQ94").Value = 40
$ws.Range("AR94").Value = 60
$ws.Range("AT94").Value = 3.15
$ws.Range("AU94").Value = 6.3
$ws.Range("AW94").Value = 5.3
$ws.Range("AX94").Value = 15
$ws.Range("AY94").Value = 18.5
$ws.Range("AZ94").Value = 65
$ws.Range("BA94").Value = 75
$ws.Range("BB94").Value = 175
$ws.Range("G94").Value = 2.15
$ws.Range("H94").Value = 3.55
$ws.Range("I94").Value = 3
$ws.Range("J94").Value = 2.7
$ws.Range("L94").Value = 3.35
$ws.Range("M94").Value = 1.04
$ws.Range("N94").Value = 8.75
$ws.Range("O94").Value = 1.19
$ws.Range("P94").Value = 4.2
$ws.Range("R94").Value = 2.25
$ws.Range("T94").Value = 3.15
$ws.Range("V94").Value = 2.42
$ws.Range("X94").Value = 12.5
$ws.Range("Y94").Value = 8.75
$ws.Range("Z94").Value = 22

# Row 95
$ws.Range("AB95").Value = 22
$ws.Range("AD95").Value = 6.8
$ws.Range("AE95").Value = 11.25
$ws.Range("AH95").Value = 11.5
$ws.Range("AI95").Value = 16.5
$ws.Range("AM95").Value = 23
$ws.Range("AN95").Value = 4.6
$ws.Range("AS95").Value = 175
$ws.Range("AT95").Value = 3
$ws.Range("AU95").Value = 6.3
$ws.Range("AV95").Value = 45
$ws.Range("AW95").Value = 4.9
$ws.Range("AX95").Value = 14
$ws.Range("AY95").Value = 18.5
$ws.Range("AZ95").Value = 55
$ws.Range("BA95").Value = 75
$ws.Range("BB95").Value = 175
$ws.Range("H95").Value = 3.4
$ws.Range("O95").Value = 1.21
$ws.Range("P95").Value = 3.95
$ws.Range("Q95").Value = 1.65
$ws.Range("R95").Value = 2.15
$ws.Range("S95").Value = 1.34
$ws.Range("T95").Value = 3
$ws.Range("U95").Value = 1.52
$ws.Range("V95").Value = 2.37
$ws.Range("W95").Value = 10.5
